$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -449.901307784372
$ws.Range("B2").Value = 8998.026

$ws.Range("A3:B4").Clear()
